$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Matthew Darby's logged time changed from "29h 50m" to "36h 33m"
$ws.Range("B5").Value = "36h 33m"

# Reflect the last-selected cell as seen in the saved workbook
$ws.Range("B5").Select()
